$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.192.24'
$ws.Range("E2").Value = '  +0.03%  '

$ws.Range("D3").Value = '1.824.48'
$ws.Range("E3").Value = '  -0.03%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.87'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.44%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5990'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.67%  '

$ws.Range("E7").Value = '  +0.18%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06941'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2762'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.80%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.47'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.22%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07600'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.54%  '

$ws.Range("D12").Value = '1.821.48'
$ws.Range("E12").Value = '  -0.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.722'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6268'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.03%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000009792'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.05%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '77.31'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.59%  '

$ws.Range("D17").Value = '28.996.60'
$ws.Range("E17").Value = '  -0.54%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.524'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -7.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.91'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -6.52%  '

$ws.Range("E20").Value = '  +0.09%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.56'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.830'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.87%  '

$ws.Range("E23").Value = '  +0.28%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '155.58'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.943'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.90%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1288'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.90%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.49'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06451'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.425'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.75%  '

$ws.Range("E30").Value = '  -1.09%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.810'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.97%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.775'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.092'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.56%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.717'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.61%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6450'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.77%  '

$ws.Range("E36").Value = '  +0.28%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.745'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.46%  '

$ws.Range("E38").Value = '  -0.53%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.582'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.12%  '

$ws.Range("D40").Value = '1.129.28'
$ws.Range("E40").Value = '  -8.37%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8914'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.57%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.003'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.20%  '

$ws.Range("D43").Value = '1.984.59'
$ws.Range("E43").Value = '  +0.73%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.42'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.42%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.07'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.43%  '

$ws.Range("E46").Value = '  -2.48%  '

$ws.Range("E47").Value = '  -0.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.446'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.53%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05501'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.43%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4526'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.347'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.21%  '
